$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = "MENAM"
$ws.Cells.Item(2, 4).Value = "WHITE"
$ws.Cells.Item(3, 2).Value = 18
$ws.Cells.Item(3, 3).Value = "NORTH SEA"
$ws.Cells.Item(3, 4).Value = "WHITE"
$ws.Cells.Item(4, 2).Value = 67
$ws.Cells.Item(4, 3).Value = "WEST AFRICA"
$ws.Cells.Item(4, 4).Value = "WHITE"
$ws.Cells.Item(5, 2).Value = 9
$ws.Cells.Item(5, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(5, 4).Value = "WHITE"
$ws.Cells.Item(6, 2).Value = 64
$ws.Cells.Item(6, 3).Value = "INDIA"
$ws.Cells.Item(6, 4).Value = "YELLOW"
$ws.Cells.Item(7, 2).Value = 265
$ws.Cells.Item(7, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(7, 4).Value = "WHITE"
$ws.Cells.Item(8, 2).Value = 275
$ws.Cells.Item(8, 3).Value = "INDIA"
$ws.Cells.Item(8, 4).Value = "WHITE"
$ws.Cells.Item(9, 2).Value = 364
$ws.Cells.Item(9, 3).Value = "NORTH SEA"
$ws.Cells.Item(9, 4).Value = "WHITE"
$ws.Cells.Item(10, 2).Value = 389
$ws.Cells.Item(10, 3).Value = "MENAM"
$ws.Cells.Item(10, 4).Value = "YELLOW"
$ws.Cells.Item(11, 2).Value = 373
$ws.Cells.Item(11, 3).Value = "WEST AFRICA"
$ws.Cells.Item(11, 4).Value = "RED"
$ws.Cells.Item(12, 2).Value = 512
$ws.Cells.Item(12, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(12, 4).Value = "WHITE"
$ws.Cells.Item(13, 2).Value = 617
$ws.Cells.Item(13, 3).Value = "WEST AFRICA"
$ws.Cells.Item(13, 4).Value = "RED"
$ws.Cells.Item(14, 2).Value = 625
$ws.Cells.Item(14, 3).Value = "NORTH SEA"
$ws.Cells.Item(14, 4).Value = "WHITE"
$ws.Cells.Item(15, 2).Value = 515
$ws.Cells.Item(15, 3).Value = "MENAM"
$ws.Cells.Item(15, 4).Value = "WHITE"
$ws.Cells.Item(16, 2).Value = 496
$ws.Cells.Item(16, 3).Value = "INDIA"
$ws.Cells.Item(16, 4).Value = "WHITE"
$ws.Cells.Item(17, 2).Value = 818
$ws.Cells.Item(17, 3).Value = "MENAM"
$ws.Cells.Item(17, 4).Value = "RED"
$ws.Cells.Item(18, 2).Value = 851
$ws.Cells.Item(18, 3).Value = "WEST AFRICA"
$ws.Cells.Item(18, 4).Value = "RED"
$ws.Cells.Item(19, 2).Value = 895
$ws.Cells.Item(19, 3).Value = "NORTH SEA"
$ws.Cells.Item(19, 4).Value = "WHITE"
$ws.Cells.Item(20, 2).Value = 849
$ws.Cells.Item(20, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(20, 4).Value = "WHITE"
$ws.Cells.Item(21, 2).Value = 907
$ws.Cells.Item(21, 3).Value = "INDIA"
$ws.Cells.Item(21, 4).Value = "RED"
$ws.Cells.Item(22, 2).Value = 1164
$ws.Cells.Item(22, 3).Value = "NORTH SEA"
$ws.Cells.Item(22, 4).Value = "WHITE"
$ws.Cells.Item(23, 2).Value = 1173
$ws.Cells.Item(23, 3).Value = "INDIA"
$ws.Cells.Item(23, 4).Value = "RED"
$ws.Cells.Item(24, 2).Value = 1266
$ws.Cells.Item(24, 3).Value = "MENAM"
$ws.Cells.Item(24, 4).Value = "RED"
$ws.Cells.Item(25, 2).Value = 1161
$ws.Cells.Item(25, 3).Value = "WEST AFRICA"
$ws.Cells.Item(25, 4).Value = "RED"
$ws.Cells.Item(26, 2).Value = 1037
$ws.Cells.Item(26, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(26, 4).Value = "WHITE"
$ws.Cells.Item(27, 2).Value = 1374
$ws.Cells.Item(27, 3).Value = "WEST AFRICA"
$ws.Cells.Item(27, 4).Value = "WHITE"
$ws.Cells.Item(28, 2).Value = 1406
$ws.Cells.Item(28, 3).Value = "NORTH SEA"
$ws.Cells.Item(28, 4).Value = "WHITE"
$ws.Cells.Item(29, 2).Value = 1305
$ws.Cells.Item(29, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(29, 4).Value = "WHITE"
$ws.Cells.Item(30, 2).Value = 1454
$ws.Cells.Item(30, 3).Value = "INDIA"
$ws.Cells.Item(30, 4).Value = "RED"
$ws.Cells.Item(31, 2).Value = 1408
$ws.Cells.Item(31, 3).Value = "MENAM"
$ws.Cells.Item(31, 4).Value = "YELLOW"
$ws.Cells.Item(32, 2).Value = 1783
$ws.Cells.Item(32, 3).Value = "INDIA"
$ws.Cells.Item(32, 4).Value = "RED"
$ws.Cells.Item(33, 2).Value = 1711
$ws.Cells.Item(33, 3).Value = "NORTH SEA"
$ws.Cells.Item(33, 4).Value = "WHITE"
$ws.Cells.Item(34, 2).Value = 1709
$ws.Cells.Item(34, 3).Value = "WEST AFRICA"
$ws.Cells.Item(34, 4).Value = "WHITE"
$ws.Cells.Item(35, 2).Value = 1603
$ws.Cells.Item(35, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(35, 4).Value = "WHITE"
$ws.Cells.Item(36, 2).Value = 1602
$ws.Cells.Item(36, 3).Value = "MENAM"
$ws.Cells.Item(36, 4).Value = "WHITE"
$ws.Cells.Item(37, 2).Value = 1996
$ws.Cells.Item(37, 3).Value = "WEST AFRICA"
$ws.Cells.Item(37, 4).Value = "WHITE"
$ws.Cells.Item(38, 2).Value = 1990
$ws.Cells.Item(38, 3).Value = "NORTH SEA"
$ws.Cells.Item(38, 4).Value = "WHITE"
$ws.Cells.Item(39, 2).Value = 1891
$ws.Cells.Item(39, 3).Value = "INDIA"
$ws.Cells.Item(39, 4).Value = "WHITE"
$ws.Cells.Item(40, 2).Value = 1889
$ws.Cells.Item(40, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(40, 4).Value = "WHITE"
$ws.Cells.Item(41, 2).Value = 1898
$ws.Cells.Item(41, 3).Value = "MENAM"
$ws.Cells.Item(41, 4).Value = "YELLOW"
$ws.Cells.Item(42, 2).Value = 2249
$ws.Cells.Item(42, 3).Value = "NORTH SEA"
$ws.Cells.Item(42, 4).Value = "WHITE"
$ws.Cells.Item(43, 2).Value = 2152
$ws.Cells.Item(43, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(43, 4).Value = "WHITE"
$ws.Cells.Item(44, 2).Value = 2139
$ws.Cells.Item(44, 3).Value = "INDIA"
$ws.Cells.Item(44, 4).Value = "WHITE"
$ws.Cells.Item(45, 2).Value = 2138
$ws.Cells.Item(45, 3).Value = "MENAM"
$ws.Cells.Item(45, 4).Value = "WHITE"
$ws.Cells.Item(46, 2).Value = 2201
$ws.Cells.Item(46, 3).Value = "WEST AFRICA"
$ws.Cells.Item(46, 4).Value = "WHITE"
$ws.Cells.Item(47, 2).Value = 2389
$ws.Cells.Item(47, 3).Value = "INDIA"
$ws.Cells.Item(47, 4).Value = "WHITE"
$ws.Cells.Item(48, 2).Value = 2449
$ws.Cells.Item(48, 3).Value = "WEST AFRICA"
$ws.Cells.Item(48, 4).Value = "WHITE"
$ws.Cells.Item(49, 2).Value = 2494
$ws.Cells.Item(49, 3).Value = "NORTH SEA"
$ws.Cells.Item(49, 4).Value = "WHITE"
$ws.Cells.Item(50, 2).Value = 2428
$ws.Cells.Item(50, 3).Value = "MENAM"
$ws.Cells.Item(50, 4).Value = "YELLOW"
$ws.Cells.Item(51, 2).Value = 2417
$ws.Cells.Item(51, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(51, 4).Value = "WHITE"
$ws.Cells.Item(52, 2).Value = 2703
$ws.Cells.Item(52, 3).Value = "WEST AFRICA"
$ws.Cells.Item(52, 4).Value = "WHITE"
$ws.Cells.Item(53, 2).Value = 2765
$ws.Cells.Item(53, 3).Value = "MENAM"
$ws.Cells.Item(53, 4).Value = "YELLOW"
$ws.Cells.Item(54, 2).Value = 2660
$ws.Cells.Item(54, 3).Value = "INDIA"
$ws.Cells.Item(54, 4).Value = "YELLOW"
$ws.Cells.Item(55, 2).Value = 2735
$ws.Cells.Item(55, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(55, 4).Value = "WHITE"
$ws.Cells.Item(56, 2).Value = 2748
$ws.Cells.Item(56, 3).Value = "NORTH SEA"
$ws.Cells.Item(56, 4).Value = "WHITE"
$ws.Cells.Item(57, 2).Value = 2903
$ws.Cells.Item(57, 3).Value = "MENAM"
$ws.Cells.Item(57, 4).Value = "WHITE"
$ws.Cells.Item(58, 2).Value = 2957
$ws.Cells.Item(58, 3).Value = "WEST AFRICA"
$ws.Cells.Item(58, 4).Value = "WHITE"
$ws.Cells.Item(59, 2).Value = 2911
$ws.Cells.Item(59, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(59, 4).Value = "WHITE"
$ws.Cells.Item(60, 2).Value = 2992
$ws.Cells.Item(60, 3).Value = "NORTH SEA"
$ws.Cells.Item(60, 4).Value = "WHITE"
$ws.Cells.Item(61, 2).Value = 2899
$ws.Cells.Item(61, 3).Value = "INDIA"
$ws.Cells.Item(61, 4).Value = "YELLOW"
$ws.Cells.Item(62, 2).Value = 3154
$ws.Cells.Item(62, 3).Value = "MENAM"
$ws.Cells.Item(62, 4).Value = "WHITE"
$ws.Cells.Item(63, 2).Value = 3337
$ws.Cells.Item(63, 3).Value = "INDIA"
$ws.Cells.Item(63, 4).Value = "YELLOW"
$ws.Cells.Item(64, 2).Value = 3215
$ws.Cells.Item(64, 3).Value = "WEST AFRICA"
$ws.Cells.Item(64, 4).Value = "WHITE"
$ws.Cells.Item(65, 2).Value = 3298
$ws.Cells.Item(65, 3).Value = "NORTH SEA"
$ws.Cells.Item(65, 4).Value = "WHITE"
$ws.Cells.Item(66, 2).Value = 3336
$ws.Cells.Item(66, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(66, 4).Value = "WHITE"
$ws.Cells.Item(67, 2).Value = 3517
$ws.Cells.Item(67, 3).Value = "INDIA"
$ws.Cells.Item(67, 4).Value = "YELLOW"
$ws.Cells.Item(68, 2).Value = 3411
$ws.Cells.Item(68, 3).Value = "MENAM"
$ws.Cells.Item(68, 4).Value = "YELLOW"
$ws.Cells.Item(69, 2).Value = 3547
$ws.Cells.Item(69, 3).Value = "NORTH SEA"
$ws.Cells.Item(69, 4).Value = "WHITE"
$ws.Cells.Item(70, 2).Value = 3462
$ws.Cells.Item(70, 3).Value = "WEST AFRICA"
$ws.Cells.Item(70, 4).Value = "WHITE"
$ws.Cells.Item(71, 2).Value = 3636
$ws.Cells.Item(71, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(71, 4).Value = "WHITE"
$ws.Cells.Item(72, 2).Value = 3737
$ws.Cells.Item(72, 3).Value = "WEST AFRICA"
$ws.Cells.Item(72, 4).Value = "WHITE"
$ws.Cells.Item(73, 2).Value = 3784
$ws.Cells.Item(73, 3).Value = "MENAM"
$ws.Cells.Item(73, 4).Value = "YELLOW"
$ws.Cells.Item(74, 2).Value = 3881
$ws.Cells.Item(74, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(74, 4).Value = "WHITE"
$ws.Cells.Item(75, 2).Value = 3827
$ws.Cells.Item(75, 3).Value = "NORTH SEA"
$ws.Cells.Item(75, 4).Value = "WHITE"
$ws.Cells.Item(76, 2).Value = 3769
$ws.Cells.Item(76, 3).Value = "INDIA"
$ws.Cells.Item(76, 4).Value = "YELLOW"
$ws.Cells.Item(77, 2).Value = 3966
$ws.Cells.Item(77, 3).Value = "MENAM"
$ws.Cells.Item(77, 4).Value = "YELLOW"
$ws.Cells.Item(78, 2).Value = 3986
$ws.Cells.Item(78, 3).Value = "INDIA"
$ws.Cells.Item(78, 4).Value = "YELLOW"
$ws.Cells.Item(79, 2).Value = 4186
$ws.Cells.Item(79, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(79, 4).Value = "WHITE"
$ws.Cells.Item(80, 2).Value = 4011
$ws.Cells.Item(80, 3).Value = "WEST AFRICA"
$ws.Cells.Item(80, 4).Value = "WHITE"
$ws.Cells.Item(81, 2).Value = 4094
$ws.Cells.Item(81, 3).Value = "NORTH SEA"
$ws.Cells.Item(81, 4).Value = "WHITE"
$ws.Cells.Item(82, 2).Value = 4311
$ws.Cells.Item(82, 3).Value = "INDIA"
$ws.Cells.Item(82, 4).Value = "YELLOW"
$ws.Cells.Item(83, 2).Value = 4397
$ws.Cells.Item(83, 3).Value = "NORTH SEA"
$ws.Cells.Item(83, 4).Value = "WHITE"
$ws.Cells.Item(84, 2).Value = 4217
$ws.Cells.Item(84, 3).Value = "MENAM"
$ws.Cells.Item(84, 4).Value = "YELLOW"
$ws.Cells.Item(85, 2).Value = 4272
$ws.Cells.Item(85, 3).Value = "WEST AFRICA"
$ws.Cells.Item(85, 4).Value = "WHITE"
$ws.Cells.Item(86, 2).Value = 4679
$ws.Cells.Item(86, 3).Value = "NORTH SEA"
$ws.Cells.Item(86, 4).Value = "WHITE"
$ws.Cells.Item(87, 2).Value = 4540
$ws.Cells.Item(87, 3).Value = "WEST AFRICA"
$ws.Cells.Item(87, 4).Value = "WHITE"
$ws.Cells.Item(88, 2).Value = 4468
$ws.Cells.Item(88, 3).Value = "MENAM"
$ws.Cells.Item(88, 4).Value = "YELLOW"
$ws.Cells.Item(89, 2).Value = 4546
$ws.Cells.Item(89, 3).Value = "INDIA"
$ws.Cells.Item(89, 4).Value = "YELLOW"
$ws.Cells.Item(90, 2).Value = 4828
$ws.Cells.Item(90, 3).Value = "WEST AFRICA"
$ws.Cells.Item(90, 4).Value = "WHITE"
$ws.Cells.Item(91, 2).Value = 4910
$ws.Cells.Item(91, 3).Value = "NORTH SEA"
$ws.Cells.Item(91, 4).Value = "WHITE"
$ws.Cells.Item(92, 2).Value = 4877
$ws.Cells.Item(92, 3).Value = "INDIA"
$ws.Cells.Item(92, 4).Value = "YELLOW"
$ws.Cells.Item(93, 2).Value = 4807
$ws.Cells.Item(93, 3).Value = "MENAM"
$ws.Cells.Item(93, 4).Value = "YELLOW"
$ws.Cells.Item(94, 2).Value = 4872
$ws.Cells.Item(94, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(94, 4).Value = "WHITE"
$ws.Cells.Item(95, 2).Value = 5192
$ws.Cells.Item(95, 3).Value = "NORTH SEA"
$ws.Cells.Item(95, 4).Value = "WHITE"
$ws.Cells.Item(96, 2).Value = 5116
$ws.Cells.Item(96, 3).Value = "WEST AFRICA"
$ws.Cells.Item(96, 4).Value = "WHITE"
$ws.Cells.Item(97, 2).Value = 5147
$ws.Cells.Item(97, 3).Value = "INDIA"
$ws.Cells.Item(97, 4).Value = "YELLOW"
$ws.Cells.Item(98, 2).Value = 5054
$ws.Cells.Item(98, 3).Value = "MENAM"
$ws.Cells.Item(98, 4).Value = "YELLOW"
$ws.Cells.Item(99, 2).Value = 5463
$ws.Cells.Item(99, 3).Value = "NORTH SEA"
$ws.Cells.Item(99, 4).Value = "WHITE"
$ws.Cells.Item(100, 2).Value = 5425
$ws.Cells.Item(100, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(100, 4).Value = "WHITE"
$ws.Cells.Item(101, 2).Value = 5393
$ws.Cells.Item(101, 3).Value = "WEST AFRICA"
$ws.Cells.Item(101, 4).Value = "WHITE"
$ws.Cells.Item(102, 2).Value = 5333
$ws.Cells.Item(102, 3).Value = "MENAM"
$ws.Cells.Item(102, 4).Value = "RED"
$ws.Cells.Item(103, 2).Value = 5391
$ws.Cells.Item(103, 3).Value = "INDIA"
$ws.Cells.Item(103, 4).Value = "YELLOW"
$ws.Cells.Item(104, 2).Value = 5611
$ws.Cells.Item(104, 3).Value = "MENAM"
$ws.Cells.Item(104, 4).Value = "YELLOW"
$ws.Cells.Item(105, 2).Value = 5742
$ws.Cells.Item(105, 3).Value = "INDIA"
$ws.Cells.Item(105, 4).Value = "YELLOW"
$ws.Cells.Item(106, 2).Value = 5674
$ws.Cells.Item(106, 3).Value = "WEST AFRICA"
$ws.Cells.Item(106, 4).Value = "WHITE"
$ws.Cells.Item(107, 2).Value = 5716
$ws.Cells.Item(107, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(107, 4).Value = "WHITE"
$ws.Cells.Item(108, 2).Value = 5766
$ws.Cells.Item(108, 3).Value = "NORTH SEA"
$ws.Cells.Item(108, 4).Value = "WHITE"
$ws.Cells.Item(109, 2).Value = 6003
$ws.Cells.Item(109, 3).Value = "INDIA"
$ws.Cells.Item(109, 4).Value = "YELLOW"
$ws.Cells.Item(110, 2).Value = 6035
$ws.Cells.Item(110, 3).Value = "SOUTH EAST ASIA"
$ws.Cells.Item(110, 4).Value = "YELLOW"
$ws.Cells.Item(111, 2).Value = 6042
$ws.Cells.Item(111, 3).Value = "MENAM"
$ws.Cells.Item(111, 4).Value = "YELLOW"
$ws.Cells.Item(112, 2).Value = 5998
$ws.Cells.Item(112, 3).Value = "WEST AFRICA"
$ws.Cells.Item(112, 4).Value = "WHITE"
$ws.Cells.Item(113, 2).Value = 6077
$ws.Cells.Item(113, 3).Value = "NORTH SEA"
$ws.Cells.Item(113, 4).Value = "WHITE"
